# Actualización automática 2025-10-16 14:30:09
#
# Updates the figures for client "F.V - AREA ANDINA S.A." / asesor
# "RIOS CARRION ANGEL BENIGNO" (PORCELANATO group) across the three
# related sheets, and nudges the "POR CUMPLIR" column width on the
# CUMPLIMIENTO MENSUAL sheet so the new, wider percentage column fits.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": PORCELANATO sales for F.V - AREA ANDINA S.A.
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M12").Value = 6935.82

# --- Sheet "VENTA MENSUAL": octubre (October) column for the same client,
#     plus the column total row.
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F12").Value = 6935.82
$wsMensual.Range("F26").Value = 23594.8

# --- Sheet "CUMPLIMIENTO MENSUAL": VENTA / POR CUMPLIR / % for the
#     PORCELANATO row and the TOTAL row, recomputed from the budget
#     (PRESUPUESTO) already stored in column C.
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D12").Value = 23374.81
$wsCumpl.Range("E12").Value = 4580.169999999998
$wsCumpl.Range("F12").Value = 0.8361590671858825

$wsCumpl.Range("D14").Value = 23594.8
$wsCumpl.Range("E14").Value = 18608.58110009469
$wsCumpl.Range("F14").Value = 0.5590736899500942

# Widen column E (POR CUMPLIR) from 22 to 23 characters. ColumnWidth uses
# Excel's character-width units, which get rounded to the nearest 1/6th
# of a character internally, so 22.17 is the safe mid-bucket value that
# lands exactly on a stored width of 23.
$wsCumpl.Columns.Item(5).ColumnWidth = 22.17
